$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value (45175) for every data
# row (2 through 357). Update it to the new serial date value (45177),
# which corresponds to the workbook being regenerated on 2023-09-08.
$ws.Range("C2:C357").Value = 45177
